$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Extra")
$ws.Range("E4").Value = "bij een item in het CMS kunnen aangeven of die zichtbaar of niet zichtbaar is"
